# "Generate Report for Handback" - fills in the Latest Target File / Latest
# Handback File / Latest Handback DateTime columns for the zh-cn and de-de
# handoff rows (now that a handback has come in), updates the Overview
# status text, and widens a few columns that now hold longer values.

$wb = $excel.ActiveWorkbook

$mdTarget1 = "66e1864f-a99b-40eb-8bff-5a3392869ca2.md"
$mdTarget2 = "8ae6cb8f-a87a-4616-933d-f94956e5ffb0.md"
$mdUrl1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4129a237e823ebb2cd9dde386170e6daab906be5/e2e/66e1864f-a99b-40eb-8bff-5a3392869ca2.md"
$mdUrl2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4129a237e823ebb2cd9dde386170e6daab906be5/e2e/8ae6cb8f-a87a-4616-933d-f94956e5ffb0.md"

# Column width achievable via ColumnWidth (engine rounds to pixel grid) that
# lands closest to the generator's 29.9777047293527 character width.
$wideWidth = 29.0833333333333
# 39.1666... rounds exactly to a stored width of 40.
$width40 = 39.1666666666667

function Update-LangSheet($sheetName, $xlfSuffix, $handbackDateTime) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Widen "Status" (C), "Latest Target File" (I) and "Latest Handback
    # File" (J) columns now that they carry longer text.
    $ws.Cells.Item(1, 3).ColumnWidth = $wideWidth
    $ws.Cells.Item(1, 9).ColumnWidth = $width40
    $ws.Cells.Item(1, 10).ColumnWidth = $width40

    # "Status" (C) moves from "In Translation" to the handed-back message.
    $ws.Range("C2").Value = "Handed back: in sync with en-US"
    $ws.Range("C3").Value = "Handed back: in sync with en-US"

    # Row 2 - 66e1864f...
    $ws.Range("I2").Value = $mdTarget1
    $ws.Range("J2").Value = "66e1864f-a99b-40eb-8bff-5a3392869ca2.a2ed611d897d7b6dcf6273f7b968c8238ea1391f.$xlfSuffix.xlf"
    $ws.Range("K2").Value = $handbackDateTime

    # Row 3 - 8ae6cb8f...
    $ws.Range("I3").Value = $mdTarget2
    $ws.Range("J3").Value = "8ae6cb8f-a87a-4616-933d-f94956e5ffb0.a8f20e30e446e058d03c3dd35a8a582f309885c1.$xlfSuffix.xlf"
    $ws.Range("K3").Value = $handbackDateTime

    # Re-create the hyperlinks so the new "Latest Target File" links (I2,
    # I3) take their place alongside the existing "Source File Name" links
    # (A2, A3), in row order.
    $ws.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range("A2"), $mdUrl1, "", "", $mdTarget1)
    $ws.Hyperlinks.Add($ws.Range("I2"), $mdUrl1, "", "", $mdTarget1)
    $ws.Hyperlinks.Add($ws.Range("A3"), $mdUrl2, "", "", $mdTarget2)
    $ws.Hyperlinks.Add($ws.Range("I3"), $mdUrl2, "", "", $mdTarget2)
}

Update-LangSheet "zh-cn" "zh-cn" "2016-08-16 12:23:53"
Update-LangSheet "de-de" "de-de" "2016-08-16 12:24:02"

# Overview sheet: the "zh-cn"/"de-de" status columns (E, F) now read
# "Handed back: in sync with en-US" and need more room.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"
$wsOverview.Cells.Item(1, 5).ColumnWidth = $wideWidth
$wsOverview.Cells.Item(1, 6).ColumnWidth = $wideWidth
